# Weekly price update for "Fruta, Agrícola del Norte S.A. de Arica - Frutilla"
# Four new daily-quality rows (date 2023-07-17 / serial 45124) are inserted
# above the existing data block (old rows 56-64 shift down to 60-68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before row 56; this pushes the former rows
# 56-64 down to 60-68 and keeps all their data/formatting intact.
$ws.Range("A56:T59").EntireRow.Insert()

# New weekly records to populate in the freshly inserted rows 56-59.
# Columns: Fecha(D), Calidad(L), Volumen(M), PrecioMin(N), PrecioMax(O), PrecioProm(P), Precio$/Kg(S)
$newRecords = @(
    @{ Fecha = 45124; Calidad = "Especial"; Volumen = 210; PrecioMin = 7000; PrecioMax = 8000; PrecioProm = 7476; PrecioKg = 2492 },
    @{ Fecha = 45124; Calidad = "Primera";  Volumen = 260; PrecioMin = 5000; PrecioMax = 6000; PrecioProm = 5462; PrecioKg = 1821 },
    @{ Fecha = 45124; Calidad = "Segunda";  Volumen = 290; PrecioMin = 4500; PrecioMax = 5000; PrecioProm = 4741; PrecioKg = 1580 },
    @{ Fecha = 45124; Calidad = "Tercera";  Volumen = 270; PrecioMin = 3500; PrecioMax = 4000; PrecioProm = 3815; PrecioKg = 1272 }
)

$r = 56
foreach ($rec in $newRecords) {
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($r, 4).Value = $rec.Fecha
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100101
    $ws.Cells.Item($r, 8).Value = "Berries"
    $ws.Cells.Item($r, 9).Value = 100112025
    $ws.Cells.Item($r, 10).Value = "Frutilla"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $rec.Calidad
    $ws.Cells.Item($r, 13).Value = $rec.Volumen
    $ws.Cells.Item($r, 14).Value = $rec.PrecioMin
    $ws.Cells.Item($r, 15).Value = $rec.PrecioMax
    $ws.Cells.Item($r, 16).Value = $rec.PrecioProm
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 3 kilos"
    $ws.Cells.Item($r, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($r, 19).Value = $rec.PrecioKg
    $ws.Cells.Item($r, 20).Value = 3
    $r++
}
